# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    populated with the Q1-2022 fund-holdings detail rows.
# 2. Insert a new summary row at the top of the "总计" sheet's data (row 2)
#    for the "2022-Q1" quarter, pushing the older quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" detail sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Helper so numeric-looking text (fund size / position figures) is kept as
# text, matching the source data, instead of being auto-coerced to a number.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - 460002 华泰柏瑞积极成长混合A
$q1.Range("A2").Value = 0
$q1.Range("A2").Font.Bold = $true
$q1.Range("A2").HorizontalAlignment = -4108
$q1.Range("A2").VerticalAlignment = -4160
$q1.Range("A2").Borders.LineStyle = 1
Set-TextValue $q1.Range("B2") "460002"
$q1.Range("C2").Value = "华泰柏瑞积极成长混合A"
Set-TextValue $q1.Range("D2") "6.11"
Set-TextValue $q1.Range("E2") "81.55"
Set-TextValue $q1.Range("F2") "2.86"
Set-TextValue $q1.Range("G2") "0.1747"
$q1.Range("H2").Value = 9

# Row 3 - 011685 创金合信先进装备股票A
$q1.Range("A3").Value = 1
$q1.Range("A3").Font.Bold = $true
$q1.Range("A3").HorizontalAlignment = -4108
$q1.Range("A3").VerticalAlignment = -4160
$q1.Range("A3").Borders.LineStyle = 1
Set-TextValue $q1.Range("B3") "011685"
$q1.Range("C3").Value = "创金合信先进装备股票A"
Set-TextValue $q1.Range("D3") "0.73"
Set-TextValue $q1.Range("E3") "92.01"
Set-TextValue $q1.Range("F3") "4.26"
Set-TextValue $q1.Range("G3") "0.0311"
$q1.Range("H3").Value = 10

# Row 4 - 011686 创金合信先进装备股票C
$q1.Range("A4").Value = 2
$q1.Range("A4").Font.Bold = $true
$q1.Range("A4").HorizontalAlignment = -4108
$q1.Range("A4").VerticalAlignment = -4160
$q1.Range("A4").Borders.LineStyle = 1
Set-TextValue $q1.Range("B4") "011686"
$q1.Range("C4").Value = "创金合信先进装备股票C"
Set-TextValue $q1.Range("D4") "0.17"
Set-TextValue $q1.Range("E4") "92.01"
Set-TextValue $q1.Range("F4") "4.26"
Set-TextValue $q1.Range("G4") "0.0072"
$q1.Range("H4").Value = 10

# Row 5 - 960030 华泰柏瑞积极成长混合H
$q1.Range("A5").Value = 3
$q1.Range("A5").Font.Bold = $true
$q1.Range("A5").HorizontalAlignment = -4108
$q1.Range("A5").VerticalAlignment = -4160
$q1.Range("A5").Borders.LineStyle = 1
Set-TextValue $q1.Range("B5") "960030"
$q1.Range("C5").Value = "华泰柏瑞积极成长混合H"
Set-TextValue $q1.Range("D5") "0.00"
Set-TextValue $q1.Range("E5") "81.55"
Set-TextValue $q1.Range("F5") "2.86"
$q1.Range("G5").Value = 0
$q1.Range("H5").Value = 9

# ---------------------------------------------------------------------------
# 2. Prepend the "2022-Q1" summary row to the "总计" sheet, shifting the
#    existing quarters down by one row and renumbering the index column.
# ---------------------------------------------------------------------------
# Re-resolve "总计" by name: after the Add() above, the worksheet collection
# shifted, so the $totalSheet handle captured earlier may no longer point at
# the right tab - looking it up again by name is the reliable way to get it.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A2").Font.Bold = $true
$totalSheet.Range("A2").HorizontalAlignment = -4108
$totalSheet.Range("A2").VerticalAlignment = -4160
$totalSheet.Range("A2").Borders.LineStyle = 1
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.21

# Renumber the existing (now shifted down) rows' index column A3:A7 -> 1..5
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
